$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wine Packages")

# Fix the "Average % Saving **" calculated column: it was dividing USD inc.
# Gratuity by the average full price, which is wrong. It should be the
# Average USD$ Saving divided by the USD inc. Gratuity actually paid.
$newFormula = "=Table2[[#This Row],[Average USD$ Saving *]]/Table2[[#This Row],[USD inc. Gratuity]]"
$ws.Range("E4").Formula = $newFormula
$ws.Range("E5").Formula = $newFormula
$ws.Range("E6").Formula = $newFormula
$ws.Range("E7").Formula = $newFormula

# Stray formatted (but empty) cell that trails along from the edit.
$ws.Range("G4").NumberFormat = "0%"

# Leave the selection on the cell that was edited.
$ws.Range("E4").Select()
